$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme
$cs.Load("Office")
